# The source added a new weekly price observation for Poroto verde at
# Vega Monumental Concepcion, inserted as the new first record (row 15).
# Every existing record from row 15 down to row 57 shifts down by one row
# (to rows 16-58), which is why the sheet dimension grows from R57 to R58.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value = 44811
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 100112031
$ws.Range("G15").Value = "Poroto verde"
$ws.Range("H15").Value = "Magnum"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 27000
$ws.Range("L15").Value = 28000
$ws.Range("M15").Value = 27500
$ws.Range("N15").Value = "$/malla 25 kilos"
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 1100
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"

# Row 16
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44265
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112031
$ws.Range("G16").Value = "Poroto verde"
$ws.Range("H16").Value = "Magnum"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 20000
$ws.Range("L16").Value = 22000
$ws.Range("M16").Value = 21000
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 840
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"

# Row 17
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44769
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112031
$ws.Range("G17").Value = "Poroto verde"
$ws.Range("H17").Value = "Magnum"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 34000
$ws.Range("L17").Value = 35000
$ws.Range("M17").Value = 34500
$ws.Range("N17").Value = "$/malla 25 kilos"
$ws.Range("O17").Value = "Perú"
$ws.Range("P17").Value = 1380
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"

# Row 18
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 44580
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100112031
$ws.Range("G18").Value = "Poroto verde"
$ws.Range("H18").Value = "Magnum"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 28000
$ws.Range("L18").Value = 30000
$ws.Range("M18").Value = 29000
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 1160
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"

# Row 19
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44708
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = "Poroto verde"
$ws.Range("H19").Value = "Magnum"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 22000
$ws.Range("M19").Value = 21000
$ws.Range("N19").Value = "$/malla 25 kilos"
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 840
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"

# Row 20
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = "Vega Monumental Concepción"
$ws.Range("C20").Value = "Bíobío"
$ws.Range("D20").Value = 44294
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 100112031
$ws.Range("G20").Value = "Poroto verde"
$ws.Range("H20").Value = "Magnum"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 24000
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = 24500
$ws.Range("N20").Value = "$/saco 25 kilos"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 980
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"

# Row 21
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44468
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 100112031
$ws.Range("G21").Value = "Poroto verde"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 31000
$ws.Range("L21").Value = 32000
$ws.Range("M21").Value = 31500
$ws.Range("N21").Value = "$/malla 25 kilos"
$ws.Range("O21").Value = "Región de Arica y Parinacota"
$ws.Range("P21").Value = 1260
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"

# Row 22
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = 44783
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 100112031
$ws.Range("G22").Value = "Poroto verde"
$ws.Range("H22").Value = "Magnum"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 38000
$ws.Range("L22").Value = 40000
$ws.Range("M22").Value = 39000
$ws.Range("N22").Value = "$/malla 25 kilos"
$ws.Range("O22").Value = "Perú"
$ws.Range("P22").Value = 1560
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"

# Row 23
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "Vega Monumental Concepción"
$ws.Range("C23").Value = "Bíobío"
$ws.Range("D23").Value = 44797
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112031
$ws.Range("G23").Value = "Poroto verde"
$ws.Range("H23").Value = "Magnum"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 44000
$ws.Range("L23").Value = 45000
$ws.Range("M23").Value = 44500
$ws.Range("N23").Value = "$/malla 25 kilos"
$ws.Range("O23").Value = "Perú"
$ws.Range("P23").Value = 1780
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"

# Row 24
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "Vega Monumental Concepción"
$ws.Range("C24").Value = "Bíobío"
$ws.Range("D24").Value = 44447
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 100112031
$ws.Range("G24").Value = "Poroto verde"
$ws.Range("H24").Value = "Magnum"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 37000
$ws.Range("L24").Value = 38000
$ws.Range("M24").Value = 37500
$ws.Range("N24").Value = "$/malla 25 kilos"
$ws.Range("O24").Value = "Perú"
$ws.Range("P24").Value = 1500
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"

# Row 25
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 44244
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100112031
$ws.Range("G25").Value = "Poroto verde"
$ws.Range("H25").Value = "Magnum"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 16000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 17000
$ws.Range("N25").Value = "$/saco 25 kilos"
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 680
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"

# Row 26
$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "Vega Monumental Concepción"
$ws.Range("C26").Value = "Bíobío"
$ws.Range("D26").Value = 44622
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 100112031
$ws.Range("G26").Value = "Poroto verde"
$ws.Range("H26").Value = "Magnum"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 220
$ws.Range("K26").Value = 24000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = 24545
$ws.Range("N26").Value = "$/malla 25 kilos"
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 982
$ws.Range("Q26").Value = 25
$ws.Range("R26").Value = "Hortaliza"

# Row 27
$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 44160
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100112031
$ws.Range("G27").Value = "Poroto verde"
$ws.Range("H27").Value = "Magnum"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 28000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = 29000
$ws.Range("N27").Value = "$/malla 25 kilos"
$ws.Range("O27").Value = "Región de O'Higgins"
$ws.Range("P27").Value = 1160
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"

# Row 28
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 44568
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112031
$ws.Range("G28").Value = "Poroto verde"
$ws.Range("H28").Value = "Magnum"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 26000
$ws.Range("M28").Value = 25500
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Región de O'Higgins"
$ws.Range("P28").Value = 1020
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"

# Row 29
$ws.Range("A29").Value = 11
$ws.Range("B29").Value = "Vega Monumental Concepción"
$ws.Range("C29").Value = "Bíobío"
$ws.Range("D29").Value = 44323
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = 100112031
$ws.Range("G29").Value = "Poroto verde"
$ws.Range("H29").Value = "Magnum"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 20000
$ws.Range("L29").Value = 22000
$ws.Range("M29").Value = 21000
$ws.Range("N29").Value = "$/malla 25 kilos"
$ws.Range("O29").Value = "Perú"
$ws.Range("P29").Value = 840
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"

# Row 30
$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 44658
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = 100112031
$ws.Range("G30").Value = "Poroto verde"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 80
$ws.Range("K30").Value = 25000
$ws.Range("L30").Value = 25000
$ws.Range("M30").Value = 25000
$ws.Range("N30").Value = "$/saco 25 kilos"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 1000
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"

# Row 31
$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = 44461
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 100112031
$ws.Range("G31").Value = "Poroto verde"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 33000
$ws.Range("L31").Value = 34000
$ws.Range("M31").Value = 33500
$ws.Range("N31").Value = "$/malla 25 kilos"
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 1340
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"

# Row 32
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44615
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112031
$ws.Range("G32").Value = "Poroto verde"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 28000
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = 29000
$ws.Range("N32").Value = "$/saco 25 kilos"
$ws.Range("O32").Value = "Región del Maule"
$ws.Range("P32").Value = 1160
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"

# Row 33
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 44574
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 100112031
$ws.Range("G33").Value = "Poroto verde"
$ws.Range("H33").Value = "Magnum"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 27000
$ws.Range("L33").Value = 28000
$ws.Range("M33").Value = 27500
$ws.Range("N33").Value = "$/saco 25 kilos"
$ws.Range("O33").Value = "Región Metropolitana"
$ws.Range("P33").Value = 1100
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"

# Row 34
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = "Vega Monumental Concepción"
$ws.Range("C34").Value = "Bíobío"
$ws.Range("D34").Value = 44706
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = 100112031
$ws.Range("G34").Value = "Poroto verde"
$ws.Range("H34").Value = "Magnum"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 100
$ws.Range("K34").Value = 22000
$ws.Range("L34").Value = 24000
$ws.Range("M34").Value = 23000
$ws.Range("N34").Value = "$/saco 25 kilos"
$ws.Range("O34").Value = "Perú"
$ws.Range("P34").Value = 920
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"

# Row 35
$ws.Range("A35").Value = 11
$ws.Range("B35").Value = "Vega Monumental Concepción"
$ws.Range("C35").Value = "Bíobío"
$ws.Range("D35").Value = 44363
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = 100112031
$ws.Range("G35").Value = "Poroto verde"
$ws.Range("H35").Value = "Magnum"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 25000
$ws.Range("L35").Value = 26000
$ws.Range("M35").Value = 25500
$ws.Range("N35").Value = "$/malla 25 kilos"
$ws.Range("O35").Value = "Perú"
$ws.Range("P35").Value = 1020
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"

# Row 36
$ws.Range("A36").Value = 11
$ws.Range("B36").Value = "Vega Monumental Concepción"
$ws.Range("C36").Value = "Bíobío"
$ws.Range("D36").Value = 44540
$ws.Range("E36").Value = 8
$ws.Range("F36").Value = 100112031
$ws.Range("G36").Value = "Poroto verde"
$ws.Range("H36").Value = "Magnum"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 170
$ws.Range("K36").Value = 21000
$ws.Range("L36").Value = 22000
$ws.Range("M36").Value = 21529
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Región Metropolitana"
$ws.Range("P36").Value = 861
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"

# Row 37
$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 44532
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 100112031
$ws.Range("G37").Value = "Poroto verde"
$ws.Range("H37").Value = "Magnum"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 250
$ws.Range("K37").Value = 33000
$ws.Range("L37").Value = 35000
$ws.Range("M37").Value = 33800
$ws.Range("N37").Value = "$/malla 25 kilos"
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 1352
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"

# Row 38
$ws.Range("A38").Value = 11
$ws.Range("B38").Value = "Vega Monumental Concepción"
$ws.Range("C38").Value = "Bíobío"
$ws.Range("D38").Value = 44384
$ws.Range("E38").Value = 8
$ws.Range("F38").Value = 100112031
$ws.Range("G38").Value = "Poroto verde"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 25000
$ws.Range("L38").Value = 26000
$ws.Range("M38").Value = 25500
$ws.Range("N38").Value = "$/malla 25 kilos"
$ws.Range("O38").Value = "Perú"
$ws.Range("P38").Value = 1020
$ws.Range("Q38").Value = 25
$ws.Range("R38").Value = "Hortaliza"

# Row 39
$ws.Range("A39").Value = 11
$ws.Range("B39").Value = "Vega Monumental Concepción"
$ws.Range("C39").Value = "Bíobío"
$ws.Range("D39").Value = 44594
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 100112031
$ws.Range("G39").Value = "Poroto verde"
$ws.Range("H39").Value = "Magnum"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 34000
$ws.Range("L39").Value = 35000
$ws.Range("M39").Value = 34500
$ws.Range("N39").Value = "$/saco 25 kilos"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 1380
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"

# Row 40
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Vega Monumental Concepción"
$ws.Range("C40").Value = "Bíobío"
$ws.Range("D40").Value = 44609
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 100112031
$ws.Range("G40").Value = "Poroto verde"
$ws.Range("H40").Value = "Magnum"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 28000
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = 29000
$ws.Range("N40").Value = "$/saco 25 kilos"
$ws.Range("O40").Value = "Región Metropolitana"
$ws.Range("P40").Value = 1160
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"

# Row 41
$ws.Range("A41").Value = 11
$ws.Range("B41").Value = "Vega Monumental Concepción"
$ws.Range("C41").Value = "Bíobío"
$ws.Range("D41").Value = 44321
$ws.Range("E41").Value = 8
$ws.Range("F41").Value = 100112031
$ws.Range("G41").Value = "Poroto verde"
$ws.Range("H41").Value = "Magnum"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 100
$ws.Range("K41").Value = 24000
$ws.Range("L41").Value = 25000
$ws.Range("M41").Value = 24500
$ws.Range("N41").Value = "$/saco 25 kilos"
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value = 980
$ws.Range("Q41").Value = 25
$ws.Range("R41").Value = "Hortaliza"

# Row 42
$ws.Range("A42").Value = 11
$ws.Range("B42").Value = "Vega Monumental Concepción"
$ws.Range("C42").Value = "Bíobío"
$ws.Range("D42").Value = 44664
$ws.Range("E42").Value = 8
$ws.Range("F42").Value = 100112031
$ws.Range("G42").Value = "Poroto verde"
$ws.Range("H42").Value = "Magnum"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 100
$ws.Range("K42").Value = 22000
$ws.Range("L42").Value = 24000
$ws.Range("M42").Value = 23000
$ws.Range("N42").Value = "$/malla 25 kilos"
$ws.Range("O42").Value = "Región Metropolitana"
$ws.Range("P42").Value = 920
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"

# Row 43
$ws.Range("A43").Value = 11
$ws.Range("B43").Value = "Vega Monumental Concepción"
$ws.Range("C43").Value = "Bíobío"
$ws.Range("D43").Value = 44651
$ws.Range("E43").Value = 8
$ws.Range("F43").Value = 100112031
$ws.Range("G43").Value = "Poroto verde"
$ws.Range("H43").Value = "Magnum"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 140
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 23000
$ws.Range("M43").Value = 21714
$ws.Range("N43").Value = "$/saco 25 kilos"
$ws.Range("O43").Value = "Región Metropolitana"
$ws.Range("P43").Value = 869
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = "Hortaliza"

# Row 44
$ws.Range("A44").Value = 11
$ws.Range("B44").Value = "Vega Monumental Concepción"
$ws.Range("C44").Value = "Bíobío"
$ws.Range("D44").Value = 44692
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 100112031
$ws.Range("G44").Value = "Poroto verde"
$ws.Range("H44").Value = "Magnum"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 100
$ws.Range("K44").Value = 25000
$ws.Range("L44").Value = 26000
$ws.Range("M44").Value = 25500
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Región Metropolitana"
$ws.Range("P44").Value = 1020
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"

# Row 45
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value = 44441
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 100112031
$ws.Range("G45").Value = "Poroto verde"
$ws.Range("H45").Value = "Magnum"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 100
$ws.Range("K45").Value = 28000
$ws.Range("L45").Value = 29000
$ws.Range("M45").Value = 28500
$ws.Range("N45").Value = "$/malla 25 kilos"
$ws.Range("O45").Value = "Perú"
$ws.Range("P45").Value = 1140
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"

# Row 46
$ws.Range("A46").Value = 11
$ws.Range("B46").Value = "Vega Monumental Concepción"
$ws.Range("C46").Value = "Bíobío"
$ws.Range("D46").Value = 44335
$ws.Range("E46").Value = 8
$ws.Range("F46").Value = 100112031
$ws.Range("G46").Value = "Poroto verde"
$ws.Range("H46").Value = "Magnum"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = 35000
$ws.Range("L46").Value = 36000
$ws.Range("M46").Value = 35500
$ws.Range("N46").Value = "$/saco 25 kilos"
$ws.Range("O46").Value = "Región Metropolitana"
$ws.Range("P46").Value = 1420
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = "Hortaliza"

# Row 47
$ws.Range("A47").Value = 11
$ws.Range("B47").Value = "Vega Monumental Concepción"
$ws.Range("C47").Value = "Bíobío"
$ws.Range("D47").Value = 44638
$ws.Range("E47").Value = 8
$ws.Range("F47").Value = 100112031
$ws.Range("G47").Value = "Poroto verde"
$ws.Range("H47").Value = "Magnum"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 180
$ws.Range("K47").Value = 23000
$ws.Range("L47").Value = 24000
$ws.Range("M47").Value = 23444
$ws.Range("N47").Value = "$/saco 25 kilos"
$ws.Range("O47").Value = "Región Metropolitana"
$ws.Range("P47").Value = 938
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"

# Row 48
$ws.Range("A48").Value = 11
$ws.Range("B48").Value = "Vega Monumental Concepción"
$ws.Range("C48").Value = "Bíobío"
$ws.Range("D48").Value = 44552
$ws.Range("E48").Value = 8
$ws.Range("F48").Value = 100112031
$ws.Range("G48").Value = "Poroto verde"
$ws.Range("H48").Value = "Magnum"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 100
$ws.Range("K48").Value = 34000
$ws.Range("L48").Value = 36000
$ws.Range("M48").Value = 35000
$ws.Range("N48").Value = "$/malla 25 kilos"
$ws.Range("O48").Value = "Región de O'Higgins"
$ws.Range("P48").Value = 1400
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"

# Row 49
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Vega Monumental Concepción"
$ws.Range("C49").Value = "Bíobío"
$ws.Range("D49").Value = 44237
$ws.Range("E49").Value = 8
$ws.Range("F49").Value = 100112031
$ws.Range("G49").Value = "Poroto verde"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 100
$ws.Range("K49").Value = 20000
$ws.Range("L49").Value = 22000
$ws.Range("M49").Value = 21000
$ws.Range("N49").Value = "$/saco 25 kilos"
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 840
$ws.Range("Q49").Value = 25
$ws.Range("R49").Value = "Hortaliza"

# Row 50
$ws.Range("A50").Value = 11
$ws.Range("B50").Value = "Vega Monumental Concepción"
$ws.Range("C50").Value = "Bíobío"
$ws.Range("D50").Value = 44279
$ws.Range("E50").Value = 8
$ws.Range("F50").Value = 100112031
$ws.Range("G50").Value = "Poroto verde"
$ws.Range("H50").Value = "Magnum"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 28000
$ws.Range("L50").Value = 30000
$ws.Range("M50").Value = 29000
$ws.Range("N50").Value = "$/saco 25 kilos"
$ws.Range("O50").Value = "Región del Maule"
$ws.Range("P50").Value = 1160
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"

# Row 51
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 44342
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112031
$ws.Range("G51").Value = "Poroto verde"
$ws.Range("H51").Value = "Magnum"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 28000
$ws.Range("L51").Value = 30000
$ws.Range("M51").Value = 29000
$ws.Range("N51").Value = "$/malla 25 kilos"
$ws.Range("O51").Value = "Región Metropolitana"
$ws.Range("P51").Value = 1160
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"

# Row 52
$ws.Range("A52").Value = 11
$ws.Range("B52").Value = "Vega Monumental Concepción"
$ws.Range("C52").Value = "Bíobío"
$ws.Range("D52").Value = 44678
$ws.Range("E52").Value = 8
$ws.Range("F52").Value = 100112031
$ws.Range("G52").Value = "Poroto verde"
$ws.Range("H52").Value = "Magnum"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 100
$ws.Range("K52").Value = 19000
$ws.Range("L52").Value = 20000
$ws.Range("M52").Value = 19500
$ws.Range("N52").Value = "$/saco 25 kilos"
$ws.Range("O52").Value = "Región Metropolitana"
$ws.Range("P52").Value = 780
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"

# Row 53
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").Value = 44629
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 100112031
$ws.Range("G53").Value = "Poroto verde"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 80
$ws.Range("K53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("M53").Value = 30000
$ws.Range("N53").Value = "$/malla 25 kilos"
$ws.Range("O53").Value = "Región Metropolitana"
$ws.Range("P53").Value = 1200
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"

# Row 54
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44272
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = 100112031
$ws.Range("G54").Value = "Poroto verde"
$ws.Range("H54").Value = "Magnum"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 100
$ws.Range("K54").Value = 22000
$ws.Range("L54").Value = 24000
$ws.Range("M54").Value = 23000
$ws.Range("N54").Value = "$/saco 25 kilos"
$ws.Range("O54").Value = "Región Metropolitana"
$ws.Range("P54").Value = 920
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"

# Row 55
$ws.Range("A55").Value = 11
$ws.Range("B55").Value = "Vega Monumental Concepción"
$ws.Range("C55").Value = "Bíobío"
$ws.Range("D55").Value = 44489
$ws.Range("E55").Value = 8
$ws.Range("F55").Value = 100112031
$ws.Range("G55").Value = "Poroto verde"
$ws.Range("H55").Value = "Magnum"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 40000
$ws.Range("L55").Value = 42000
$ws.Range("M55").Value = 41000
$ws.Range("N55").Value = "$/malla 25 kilos"
$ws.Range("O55").Value = "Perú"
$ws.Range("P55").Value = 1640
$ws.Range("Q55").Value = 25
$ws.Range("R55").Value = "Hortaliza"

# Row 56
$ws.Range("A56").Value = 11
$ws.Range("B56").Value = "Vega Monumental Concepción"
$ws.Range("C56").Value = "Bíobío"
$ws.Range("D56").Value = 44636
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 100112031
$ws.Range("G56").Value = "Poroto verde"
$ws.Range("H56").Value = "Magnum"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 180
$ws.Range("K56").Value = 22000
$ws.Range("L56").Value = 23000
$ws.Range("M56").Value = 22444
$ws.Range("N56").Value = "$/saco 25 kilos"
$ws.Range("O56").Value = "Región del Maule"
$ws.Range("P56").Value = 898
$ws.Range("Q56").Value = 25
$ws.Range("R56").Value = "Hortaliza"

# Row 57
$ws.Range("A57").Value = 11
$ws.Range("B57").Value = "Vega Monumental Concepción"
$ws.Range("C57").Value = "Bíobío"
$ws.Range("D57").Value = 44435
$ws.Range("E57").Value = 8
$ws.Range("F57").Value = 100112031
$ws.Range("G57").Value = "Poroto verde"
$ws.Range("H57").Value = "Magnum"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 100
$ws.Range("K57").Value = 25000
$ws.Range("L57").Value = 26000
$ws.Range("M57").Value = 25500
$ws.Range("N57").Value = "$/malla 25 kilos"
$ws.Range("O57").Value = "Perú"
$ws.Range("P57").Value = 1020
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"

# Row 58
$ws.Range("A58").Value = 11
$ws.Range("B58").Value = "Vega Monumental Concepción"
$ws.Range("C58").Value = "Bíobío"
$ws.Range("D58").Value = 44510
$ws.Range("E58").Value = 8
$ws.Range("F58").Value = 100112031
$ws.Range("G58").Value = "Poroto verde"
$ws.Range("H58").Value = "Magnum"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 100
$ws.Range("K58").Value = 35000
$ws.Range("L58").Value = 36000
$ws.Range("M58").Value = 35500
$ws.Range("N58").Value = "$/malla 25 kilos"
$ws.Range("O58").Value = "Perú"
$ws.Range("P58").Value = 1420
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"

# Row 58 is brand new in this sheet; give its date cell (column D) the same
# date number format already used by the rest of column D.
$ws.Range("D58").NumberFormat = $ws.Range("D57").NumberFormat
